# Auto commit at 2026-02-09 10:20:58.69
#
# Updates the "Metrics" sheet's raw metric values (B2:B13). The "today"
# sheet's B11:B22 (and their dependent E/F columns) are live formulas that
# reference Metrics!B2:B13, and A1 there recomputes TODAY()-1, so they all
# ripple through automatically on recalculation - no need to touch them
# directly. Finally, restore/update the two sheets' saved cell selections.

$wb = $excel.ActiveWorkbook

$wsMetrics = $wb.Worksheets.Item("Metrics")
$wsToday   = $wb.Worksheets.Item("today")

# --- Metrics!B2:B13 -> new values -----------------------------------------
$wsMetrics.Range("B2").Value  = 124300.93999999999
$wsMetrics.Range("B3").Value  = 112244.35999999999
$wsMetrics.Range("B4").Value  = 42595.73
$wsMetrics.Range("B5").Value  = 5012
$wsMetrics.Range("B6").Value  = 704936.71
$wsMetrics.Range("B7").Value  = 565038.81000000006
$wsMetrics.Range("B8").Value  = 207030.29
$wsMetrics.Range("B9").Value  = 28488
$wsMetrics.Range("B10").Value = 34806188.43
$wsMetrics.Range("B11").Value = 32611031.599999998
$wsMetrics.Range("B12").Value = 12152844.149999999
$wsMetrics.Range("B13").Value = 1346395

# --- Selections -------------------------------------------------------------
# Set the Metrics selection first (this does not disturb which sheet is the
# active tab), then re-activate "today" (the originally active tab) before
# updating its own selection so tabSelected stays put on "today".
$wsMetrics.Range("E20").Select()

$wsToday.Activate()
$wsToday.Range("E7").Select()
